$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (column D) and Volume(1h) (column E) cells with the latest
# scraped crypto values. Price strings that look like plain decimal numbers
# (e.g. "95.44") are written with the cell pre-formatted as Text ("@") so Excel
# keeps them as literal text instead of auto-converting them to numeric values
# (values that already use "." as a thousands separator, e.g. "43.577.90", are
# not valid numbers so Excel stores them as text without any extra formatting).

$ws.Range("D2").Value = "43.577.90"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").Value = "2.280.17"
$ws.Range("E3").Value = "  +1.60%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "95.44"
$ws.Range("E5").Value = "  -3.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.16"
$ws.Range("E6").Value = "  -1.50%  "

$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -4.55%  "

$ws.Range("E10").Value = "  -8.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.71"
$ws.Range("E12").Value = "  -8.05%  "

$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("D14").Value = "2.619.94"
$ws.Range("E14").Value = "  +1.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.16"
$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.847"
$ws.Range("E16").Value = "  +2.59%  "

$ws.Range("D17").Value = "2.289.44"
$ws.Range("E17").Value = "  +2.04%  "

$ws.Range("D18").Value = "43.590.50"
$ws.Range("E18").Value = "  -0.85%  "

$ws.Range("E19").Value = "  +2.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("E21").Value = "  +1.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.36"
$ws.Range("E22").Value = "  +1.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.92"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("E24").Value = "  -2.78%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.37"
$ws.Range("E26").Value = "  -1.14%  "

$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.45"
$ws.Range("E28").Value = "  -2.49%  "

$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.52"
$ws.Range("E30").Value = "  -3.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "176.25"
$ws.Range("E31").Value = "  +2.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.04"
$ws.Range("E32").Value = "  +4.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0885"
$ws.Range("E33").Value = "  -3.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.36"
$ws.Range("E34").Value = "  -4.15%  "

$ws.Range("E35").Value = "  +0.62%  "

$ws.Range("E36").Value = "  -3.77%  "

$ws.Range("E37").Value = "  +0.89%  "

$ws.Range("E38").Value = "  +2.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.27"
$ws.Range("E39").Value = "  -11.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.37"
$ws.Range("E40").Value = "  +8.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.237"
$ws.Range("E41").Value = "  -5.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.36"
$ws.Range("E42").Value = "  +17.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.90"
$ws.Range("E43").Value = "  -5.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.77"
$ws.Range("E44").Value = "  +0.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.84"
$ws.Range("E45").Value = "  +4.12%  "

$ws.Range("E46").Value = "  -3.97%  "

$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.40"
$ws.Range("E48").Value = "  -2.11%  "

$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.54"
$ws.Range("E50").Value = "  +7.72%  "

$ws.Range("D51").Value = "2.501.28"
$ws.Range("E51").Value = "  +1.69%  "
